$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 66.875
$ws.Range("I2").Value = 63.863636
$ws.Range("K2").Value = 63.863636
$ws.Range("M2").Value = 49.136364

$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

$ws.Range("H9").Value = 499.94736
$ws.Range("J9").Value = 629.6667
$ws.Range("L9").Value = 629.6667
$ws.Range("N9").Value = -967.6667

$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()

$ws.Range("H28").Value = 976.1905
$ws.Range("I28").Value = 578.6316
$ws.Range("K28").Value = 578.6316
$ws.Range("M28").Value = -93.63160000000005

$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()

$ws.Range("H74").Value = 60014376
$ws.Range("I74").Value = 93761864
$ws.Range("K74").Value = 93761864
$ws.Range("M74").Value = -93760928

$ws.Range("H77").Value = 60014376
$ws.Range("I77").Value = 93761864
$ws.Range("K77").Value = 468809320
$ws.Range("M77").Value = -468804640

$ws.Range("H138").Value = 2448.88
$ws.Range("I138").Value = 2002.7567
$ws.Range("K138").Value = 6008.2701
$ws.Range("M138").Value = -868.2700999999997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 16395453
$ws.Range("I2").Value = 923.4792
$ws.Range("J2").Value = 76929100
$ws.Range("K2").Value = 923.4792
$ws.Range("L2").Value = 76929100
$ws.Range("M2").Value = -810.4792
$ws.Range("N2").Value = -76929326

$ws.Range("H32").Value = 2789545.8
$ws.Range("I32").Value = 3395917.8
$ws.Range("J32").Value = 37550.23
$ws.Range("K32").Value = 3395917.8
$ws.Range("L32").Value = 37550.23
$ws.Range("M32").Value = -3395630.8
$ws.Range("N32").Value = -38124.23

$ws.Range("H45").Value = 4455.1113
$ws.Range("J45").Value = 4788
$ws.Range("L45").Value = 4788
$ws.Range("N45").Value = -5542

$ws.Range("H61").Value = 27030074
$ws.Range("I61").Value = 2145.9583
$ws.Range("J61").Value = 76927784
$ws.Range("K61").Value = 2145.9583
$ws.Range("L61").Value = 76927784
$ws.Range("M61").Value = -1933.9583
$ws.Range("N61").Value = -76928208

$ws.Range("H63").Value = 1611.5883
$ws.Range("I63").Value = 1638.3077
$ws.Range("J63").Value = 1524.75
$ws.Range("K63").Value = 1638.3077
$ws.Range("L63").Value = 1524.75
$ws.Range("M63").Value = -952.3077000000001
$ws.Range("N63").Value = -2896.75

$ws.Range("H66").Value = 1611.5883
$ws.Range("I66").Value = 1638.3077
$ws.Range("J66").Value = 1524.75
$ws.Range("K66").Value = 8191.538500000001
$ws.Range("L66").Value = 7623.75
$ws.Range("M66").Value = -4759.538500000001
$ws.Range("N66").Value = -14487.75

$ws.Range("H74").Value = 37754.656
$ws.Range("I74").Value = 64321.438
$ws.Range("J74").Value = 5057.077
$ws.Range("K74").Value = 64321.438
$ws.Range("L74").Value = 5057.077
$ws.Range("M74").Value = -63447.438
$ws.Range("N74").Value = -6805.077

$ws.Range("H77").Value = 37754.656
$ws.Range("I77").Value = 64321.438
$ws.Range("J77").Value = 5057.077
$ws.Range("K77").Value = 321607.19
$ws.Range("L77").Value = 25285.385
$ws.Range("M77").Value = -317239.19
$ws.Range("N77").Value = -34021.385

$ws.Range("H116").Value = 16395453
$ws.Range("I116").Value = 923.4792
$ws.Range("J116").Value = 76929100
$ws.Range("K116").Value = 923.4792
$ws.Range("L116").Value = 76929100
$ws.Range("M116").Value = 1370.5208
$ws.Range("N116").Value = -76933688

$ws.Range("H132").Value = 2844.8215
$ws.Range("I132").Value = 1188.4474
$ws.Range("J132").Value = 6341.6113
$ws.Range("K132").Value = 3565.3422
$ws.Range("L132").Value = 19024.8339
$ws.Range("M132").Value = -1035.3422
$ws.Range("N132").Value = -24084.8339

$ws.Range("H135").Value = 88446.38
$ws.Range("J135").Value = 88446.38
$ws.Range("L135").Value = 88446.38
$ws.Range("N135").Value = -98586.38

$ws.Range("H136").Value = 27030074
$ws.Range("I136").Value = 2145.9583
$ws.Range("J136").Value = 76927784
$ws.Range("K136").Value = 6437.874899999999
$ws.Range("L136").Value = 230783352
$ws.Range("M136").Value = -3887.874899999999
$ws.Range("N136").Value = -230788452

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 16395453
$ws.Range("I3").Value = 923.4792
$ws.Range("J3").Value = 76929100
$ws.Range("K3").Value = 923.4792
$ws.Range("L3").Value = 76929100
$ws.Range("M3").Value = -809.4792
$ws.Range("N3").Value = -76929328

$ws.Range("H5").Value = 1031.125
$ws.Range("I5").Value = 500
$ws.Range("J5").Value = 1349.8
$ws.Range("K5").Value = 500
$ws.Range("L5").Value = 1349.8
$ws.Range("M5").Value = -387
$ws.Range("N5").Value = -1575.8

$ws.Range("H113").Value = 4965
$ws.Range("I113").Value = 4965
$ws.Range("K113").Value = 4965
$ws.Range("M113").Value = -2795

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5227.245
$ws.Range("I31").Value = 1891.25
$ws.Range("K31").Value = 1891.25
$ws.Range("M31").Value = -1596.25

$ws.Range("H34").Value = 5227.245
$ws.Range("I34").Value = 1891.25
$ws.Range("K34").Value = 1891.25
$ws.Range("M34").Value = -1689.25

$ws.Range("H58").Value = 4312.185
$ws.Range("J58").Value = 5953.0415
$ws.Range("L58").Value = 5953.0415
$ws.Range("N58").Value = -6359.0415

$ws.Range("H132").Value = 4175.6045
$ws.Range("I132").Value = 2772.1482
$ws.Range("J132").Value = 6543.9375
$ws.Range("K132").Value = 8316.444600000001
$ws.Range("L132").Value = 19631.8125
$ws.Range("M132").Value = -5786.444600000001
$ws.Range("N132").Value = -24691.8125

$ws.Range("H136").Value = 4312.185
$ws.Range("J136").Value = 5953.0415
$ws.Range("L136").Value = 17859.1245
$ws.Range("N136").Value = -22959.1245

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 1004
$ws.Range("I51").Value = 1004
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 3012
$ws.Range("L51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -2552

$ws.Range("H122").Value = 3144303
$ws.Range("I122").Value = 4042418.2
$ws.Range("J122").Value = 900
$ws.Range("K122").Value = 36381763.8
$ws.Range("L122").Value = 8100
$ws.Range("M122").Value = -36379313.8
$ws.Range("N122").Value = -13000

$ws.Range("H129").Value = 168828.5
$ws.Range("J129").Value = 202288.3
$ws.Range("L129").Value = 606864.8999999999
$ws.Range("N129").Value = -616864.8999999999

$ws.Range("H139").Value = 42496.555
$ws.Range("I139").Value = 61730.53
$ws.Range("J139").Value = 9798.799999999999
$ws.Range("K139").Value = 185191.59
$ws.Range("L139").Value = 29396.4
$ws.Range("M139").Value = -180051.59
$ws.Range("N139").Value = -39676.39999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 774.5
$ws.Range("J3").Value = 774.5
$ws.Range("L3").Value = 774.5
$ws.Range("N3").Value = -1006.5

$ws.Range("H12").Value = 1007500
$ws.Range("I12").Value = 1007500
$ws.Range("K12").Value = 1007500
$ws.Range("M12").Value = -1007360

$ws.Range("J102").Value = 9000
$ws.Range("L102").Value = 9000
$ws.Range("N102").Value = -12244

$ws.Range("H118").Value = 39000
$ws.Range("J118").Value = 39000
$ws.Range("L118").Value = 39000
$ws.Range("N118").Value = -42314

$ws.Range("H126").Value = 4104.4443
$ws.Range("J126").Value = 4158.6924
$ws.Range("L126").Value = 12476.0772
$ws.Range("N126").Value = -17416.0772

$ws.Range("H132").Value = 2267
$ws.Range("I132").Value = 1210.4667
$ws.Range("J132").Value = 4908.3335
$ws.Range("K132").Value = 3631.4001
$ws.Range("L132").Value = 14725.0005
$ws.Range("M132").Value = -1101.4001
$ws.Range("N132").Value = -19785.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 239.26666
$ws.Range("I55").Value = 255.93333
$ws.Range("K55").Value = 255.93333
$ws.Range("M55").Value = -82.93333000000001

$ws.Range("H132").Value = 8776166
$ws.Range("I132").Value = 15628009
$ws.Range("K132").Value = 46884027
$ws.Range("M132").Value = -46881497

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 1336
$ws.Range("I7").Value = 1504
$ws.Range("J7").Value = 1000
$ws.Range("K7").Value = 1504
$ws.Range("L7").Value = 1000
$ws.Range("M7").Value = -1391
$ws.Range("N7").Value = -1226

$ws.Range("H51").Value = 15487.5
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H52").Value = 12900
$ws.Range("I52").Value = 12900
$ws.Range("K52").Value = 12900
$ws.Range("M52").Value = -12674

$ws.Range("H55").Value = 6459.4
$ws.Range("J55").Value = 9000
$ws.Range("L55").Value = 9000
$ws.Range("N55").Value = -9554

$ws.Range("H113").Value = 8342.815000000001
$ws.Range("J113").Value = 2029
$ws.Range("L113").Value = 6087
$ws.Range("N113").Value = -10427

$ws.Range("H117").Value = 33183.5
$ws.Range("I117").Value = 10000
$ws.Range("J117").Value = 40911.332
$ws.Range("K117").Value = 10000
$ws.Range("L117").Value = 40911.332
$ws.Range("M117").Value = -5411
$ws.Range("N117").Value = -50089.332

$ws.Range("H132").Value = 4219.102
$ws.Range("I132").Value = 3818.2632
$ws.Range("K132").Value = 11454.7896
$ws.Range("M132").Value = -8924.7896
